$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that sits right after "HVAC.HVAC Controls".
#    Word auto-renumbers every bookmark id that follows (ids 4..12 shift down to 3..11).
$d.Bookmarks("_GoBack").Delete()

# 2. Replace the lone "n/a" paragraph text with the new sentence.
[void]$d.Content.Find.Execute("n/a", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Works with Schedule Ruleset.", 2)

# 3. Re-insert a "_GoBack" bookmark, now between "Define arguments (" and "zone".
$r = $d.Content
[void]$r.Find.Execute("Define arguments (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
